$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Internships 2022")

# Insert 1 row at row 3
$ws.Rows.Item(3).Insert()

# Insert 2 rows at row 6
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# Fix formatting on the newly inserted rows by copying format from row 4 (A4:D4), a known "s=3" row
$ws.Range("A4:D4").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)
$ws.Range("A6:D7").PasteSpecial(-4122)

Write-Output "done"
